# Daily attendance processing - 2026-01-22 07:17:51
# Swap the order of "System" and the email address in the "Last Edited By"
# (column G) entries: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Only touch the individual cells that actually contain the exact text,
# leaving every other (including blank) cell in column G untouched.
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value2 = $newText
    }
}
